$d = $word.ActiveDocument

# Locate the paragraph holding the site footer/copyright text
# ("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
# pages. Original theme under Creative Commons Attribution").
$rng = $d.Content
$rng.Find.Execute("Contact: luizeleno@usp.br", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Remove the footer paragraph together with the two empty spacer /
    # forced-page-break paragraphs that immediately precede it.
    $startIndex = $targetIndex - 2
    if ($startIndex -lt 1) { $startIndex = 1 }

    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($targetIndex)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
